{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// Stable references to the original paragraphs we need to touch / anchor off of.\nconst pProblem = paras.items[1];    // \"[EM] Consolidation ... 1. Customer Problem...\"\nconst pResearch = paras.items[2];   // \"2. Customer Research...\"\nconst pSolution = paras.items[3];   // \"3. Our Solution...\"\nconst pMetrics = paras.items[7];    // \"4. Product Metrics...\"\nconst pPrototype = paras.items[13]; // \"Appendix: Quick prototype\" (Heading2)\n\n// ---- Step 1: insert the new paragraph structure first (while surrounding\n// formatting is still plain, so new paragraphs don't inherit bold/spacing) ----\n\n// Empty paragraph after the problem statement.\npProblem.insertParagraph(\"\", Word.InsertLocation.after);\n\n// After the research paragraph: empty, two new body paragraphs, empty.\nlet afterResearch = pResearch.insertParagraph(\"\", Word.InsertLocation.after);\nlet pFeedback = afterResearch.insertParagraph(\n  \"This capability was requested as feedback from an enterprise-level accounting firm, reflecting needs observed in large multi-entity audit workflows.\",\n  Word.InsertLocation.after\n);\nlet pParity = pFeedback.insertParagraph(\n  \"We are also building this to achieve competitive parity with Wolters Kluwer ProSystem fx Engagement, which offers similar functionality.\",\n  Word.InsertLocation.after\n);\npParity.insertParagraph(\"\", Word.InsertLocation.after);\n\n// Empty paragraph after the solution paragraph.\npSolution.insertParagraph(\"\", Word.InsertLocation.after);\n\n// Empty paragraphs before/after the metrics paragraph.\npMetrics.insertParagraph(\"\", Word.InsertLocation.before);\npMetrics.insertParagraph(\"\", Word.InsertLocation.after);\n\n// Empty paragraphs before/after the prototype heading; these must NOT inherit\n// the Heading2 style, so reset them to Normal.\nconst pBeforeProto = pPrototype.insertParagraph(\"\", Word.InsertLocation.before);\npBeforeProto.style = \"Normal\";\nconst pAfterProto = pPrototype.insertParagraph(\"\", Word.InsertLocation.after);\npAfterProto.style = \"Normal\";\n\nawait context.sync();\n\n// ---- Step 2: bold the run text (not the paragraph mark) and add\n// before/after spacing on the 4 \"section header\" paragraphs ----\nasync function boldAndSpace(p) {\n  p.load(\"text\");\n  await context.sync();\n  const sr = body.search(p.text, { matchCase: true });\n  sr.load(\"items\");\n  await context.sync();\n  sr.items[0].font.bold = true;\n  p.spaceBefore = 6;\n  p.spaceAfter = 6;\n  await context.sync();\n}\n\nawait boldAndSpace(pResearch);\nawait boldAndSpace(pSolution);\nawait boldAndSpace(pMetrics);\nawait boldAndSpace(pPrototype);\n", "ps1": "$d = $word.ActiveDocument\n\n# Work from the bottom of the document upward so that paragraph indices for\n# not-yet-processed (earlier) sections stay stable while we edit later ones.\n\n# --- Section: Appendix: Quick prototype heading (originally paragraph 14) ---\n$p14 = $d.Paragraphs.Item(14)\n$p14.Range.InsertParagraphAfter()\n$protoAfter = $d.Paragraphs.Item(15)\n$protoAfter.Style = \"Normal\"\n\n$p14.Range.InsertParagraphBefore()\n$protoBefore = $d.Paragraphs.Item(14)\n$protoBefore.Style = \"Normal\"\n\n$proto = $d.Paragraphs.Item(15)\n$rngProto = $proto.Range.Duplicate\n[void]$rngProto.MoveEnd(1, -1)\n$rngProto.Font.Bold = 1\n$proto.SpaceBefore = 6\n$proto.SpaceAfter = 6\n\n# --- Section: 4. Product Metrics (originally paragraph 8) ---\n$p8 = $d.Paragraphs.Item(8)\n$p8.Range.InsertParagraphAfter()\n$p8.Range.InsertParagraphBefore()\n$metrics = $d.Paragraphs.Item(9)\n$rngMetrics = $metrics.Range.Duplicate\n[void]$rngMetrics.MoveEnd(1, -1)\n$rngMetrics.Font.Bold = 1\n$metrics.SpaceBefore = 6\n$metrics.SpaceAfter = 6\n\n# --- Section: 3. Our Solution (originally paragraph 4) ---\n$p4 = $d.Paragraphs.Item(4)\n$p4.Range.InsertParagraphAfter()\n$rngSol = $p4.Range.Duplicate\n[void]$rngSol.MoveEnd(1, -1)\n$rngSol.Font.Bold = 1\n$p4.SpaceBefore = 6\n$p4.SpaceAfter = 6\n\n# --- Section: 2. Customer Research (originally paragraph 3) ---\n$p3 = $d.Paragraphs.Item(3)\n\n$p3.Range.InsertParagraphAfter()\n$emptyAfterResearch = $d.Paragraphs.Item(4)\n$emptyAfterResearch.Range.InsertParagraphAfter()\n$txt1 = $d.Paragraphs.Item(5)\n$txt1.Range.InsertBefore(\"This capability was requested as feedback from an enterprise-level accounting firm, reflecting needs observed in large multi-entity audit workflows.\")\n\n$txt1.Range.InsertParagraphAfter()\n$txt2 = $d.Paragraphs.Item(6)\n$txt2.Range.InsertBefore(\"We are also building this to achieve competitive parity with Wolters Kluwer ProSystem fx Engagement, which offers similar functionality.\")\n\n$txt2.Range.InsertParagraphAfter()\n\n$rngRes = $p3.Range.Duplicate\n[void]$rngRes.MoveEnd(1, -1)\n$rngRes.Font.Bold = 1\n$p3.SpaceBefore = 6\n$p3.SpaceAfter = 6\n\n# --- Empty paragraph after the customer-problem statement (originally paragraph 2) ---\n$p2 = $d.Paragraphs.Item(2)\n$p2.Range.InsertParagraphAfter()\n"}
